$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original General-text display by
# writing through a Text number format, then resetting the style back to Normal
# so no stray per-cell style index lingers on the written cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.872.89'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.891.86'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '0.7757'
$ws.Range("E5").Value = '  -1.68%  '
$ws.Range("D6").Value = '244.13'
$ws.Range("E6").Value = '  +0.75%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").Value = '0.3140'
$ws.Range("E8").Value = '  -1.27%  '
$ws.Range("D9").Value = '0.07389'
$ws.Range("E9").Value = '  +4.55%  '
$ws.Range("D10").Value = '25.30'
$ws.Range("E10").Value = '  -1.47%  '
$ws.Range("D11").Value = '0.08144'
$ws.Range("E11").Value = '  +1.06%  '
$ws.Range("D12").Value = '0.7659'
$ws.Range("E12").Value = '  -0.60%  '
$ws.Range("D13").Value = '5.475'
$ws.Range("E13").Value = '  +3.20%  '
$ws.Range("D14").Value = '1.853.00'
$ws.Range("E14").Value = '  -1.89%  '
$ws.Range("D16").Value = '6.218'
$ws.Range("E16").Value = '  +5.10%  '
$ws.Range("D17").Value = '29.834.33'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("E18").Value = '  +0.76%  '
$ws.Range("D19").Value = '245.27'
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").Value = '0.000007862'
$ws.Range("E20").Value = '  +1.96%  '
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").Value = '8.131'
$ws.Range("E22").Value = '  +0.90%  '
$ws.Range("D23").Value = '2.111.69'
$ws.Range("E23").Value = '  -1.36%  '
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("D26").Value = '9.429'
$ws.Range("E26").Value = '  +1.47%  '
$ws.Range("D27").Value = '162.18'
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("D28").Value = '18.79'
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("D29").Value = '2.038'
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").Value = '1.456'
$ws.Range("E30").Value = '  +5.70%  '
$ws.Range("D31").Value = '1.545'
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("D32").Value = '4.493'
$ws.Range("E32").Value = '  +1.72%  '
$ws.Range("D33").Value = '0.05599'
$ws.Range("E33").Value = '  -0.49%  '
$ws.Range("D34").Value = '4.097'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").Value = '1.251'
$ws.Range("E35").Value = '  -1.11%  '
$ws.Range("D36").Value = '0.7576'
$ws.Range("E36").Value = '  +2.90%  '
$ws.Range("D37").Value = '0.9982'
$ws.Range("E37").Value = '  -0.45%  '
$ws.Range("D38").Value = '2.648'
$ws.Range("E38").Value = '  -2.16%  '
$ws.Range("D39").Value = '0.01935'
$ws.Range("E39").Value = '  +0.53%  '
$ws.Range("D40").Value = '2.791'
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("D41").Value = '1.148.89'
$ws.Range("E41").Value = '  +12.49%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.4458'
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '74.09'
$ws.Range("E43").Value = '  +2.67%  '
$ws.Range("D44").Value = '5.962'
$ws.Range("E44").Value = '  +1.66%  '
$ws.Range("D45").Value = '0.8534'
$ws.Range("E45").Value = '  +1.02%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '1.903'
$ws.Range("E47").Value = '  +1.33%  '
$ws.Range("E48").Value = '  +6.16%  '
$ws.Range("D49").Value = '101.97'
$ws.Range("E49").Value = '  -0.33%  '
$ws.Range("D50").Value = '9.867'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("D51").Value = '7.515'
$ws.Range("E51").Value = '  +0.66%  '

$ws.Range("D2:E51").Style = "Normal"
